$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The template currently has a single (empty) data row at row 4, a totals
# row at row 5, and a footer row at row 6. We need 14 data rows, so insert
# 13 more rows right after row 4 (pushing the totals/footer rows down to
# row 18 / row 19) and fill them all in.
# ---------------------------------------------------------------------------

$ws.Range("A5:A17").EntireRow.Insert()

# Copy row 4's cell formatting (font/fill/border/alignment) into the new
# rows so every data row looks the same as the original template row.
$ws.Range("A4:N4").Copy()
for ($r = 5; $r -le 17; $r++) {
    $ws.Range("A" + $r + ":N" + $r).PasteSpecial(-4122)
}

# Re-create the merged cell layout (B:G, H:K, L:M) for every data row.
for ($r = 4; $r -le 17; $r++) {
    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
}

# Row heights, as captured by the source report for this data set.
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 24.75
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 24.75
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 25.5
$ws.Rows.Item(19).RowHeight = 17.25

# Column B (item name) and column H (transaction ratio, e.g. "12:0") must
# stay text - otherwise values like "0:0" could be reinterpreted.
$ws.Range("B4:B17").NumberFormat = "@"
$ws.Range("H4:H17").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------

$data = @(
    @(1,  "ANGIOFOX (EFFOX) 25MG LONG 30 CAPS.",        "0:0",    114,    1),
    @(2,  "AUGMENTIN 457MG/5ML SUSP. 70 ML",             "1:0",    137,    1),
    @(3,  "BLOKATENS 10/160MG 28 F.C.TABS.",             "0:0",    160,    1),
    @(4,  "COLOVATIL 30 F.C. TABS",                      "0:0",    63,     1),
    @(5,  "GAVISCON LIQUID 24 SACHETS 10 ML",            "0:9",    12,     0.04),
    @(6,  "GINKGO BILOBA 30 CAPS.",                      "0:0",    186,    1),
    @(7,  "MILGA ADVANCE 30 F.C. TABS",                  "0:0",    136.5,  1),
    @(8,  "PERLOC 40MG 14 F.C.TAB.",                     "0:0",    68.25,  1),
    @(9,  "RHINEX 0.05% INFANTILE NASAL DROPS 10 ML",    "2:0",    18,     1),
    @(10, "RIVO 320MG 20*10 TABS",                       "1:2",    14.1,   0.1),
    @(11, "VASTAREL MR 35MG 30 F.C.TAB.",                "2:0",    175,    1),
    @(12, "WATER FOR INJECTION AMP. 5 ML",               "7816:0", 2.5,    1),
    @(13, "سويت كوكو",                                    "22:0",   25,     1),
    @(14, "مرطب شفاه لونا جوز هند ابيض",                  "3:0",    20,     1)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("H" + $r).Value = $row[2]
    $ws.Range("L" + $r).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $r = $r + 1
}

# Totals row (previously row 5, now row 18).
$ws.Range("K18").Value = 1131.35

Write-Output "edit complete"
